# lsh_coding.xlsx: add a new "clinical_assessment_categories" sheet right
# before "lsh_sheet_names" (which becomes the new active/selected tab),
# and tweak a couple of cosmetic view properties on existing sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new sheet just before "lsh_sheet_names" so the tab order
#    becomes: lsh_covid_groups, lsh_services, lsh_unit_categories,
#    lsh_text_out_categories, clinical_assessment_categories,
#    lsh_sheet_names
# ---------------------------------------------------------------------
$namesSheetBefore = $wb.Worksheets.Item("lsh_sheet_names")
$ws = $wb.Worksheets.Add($namesSheetBefore)
$ws.Name = "clinical_assessment_categories"

# ---------------------------------------------------------------------
# 2. Populate the new sheet with the Icelandic/English clinical
#    assessment (triage colour) category lookup table. Values are
#    written header-row, then column A, column B, column C (six base
#    colours), then the "Gulur/yellow" row that was appended later, and
#    finally column D - matching how the sheet was actually authored.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "clinical_assessment_category_raw"
$ws.Range("B1").Value = "clinical_assessment_category_all"
$ws.Range("C1").Value = "clinical_assessment_category_simple"

$ws.Range("A2").Value = "Blár"
$ws.Range("A3").Value = "Brúnn"
$ws.Range("A4").Value = "Grænn"
$ws.Range("A6").Value = "Hvítur"
$ws.Range("A7").Value = "Óskilgreint"
$ws.Range("A8").Value = "Rauður"

$ws.Range("B2").Value = "blue"
$ws.Range("B3").Value = "brown"
$ws.Range("B4").Value = "green"
$ws.Range("B6").Value = "white"
$ws.Range("B7").Value = "unknown"
$ws.Range("B8").Value = "red"

$ws.Range("C2").Value = "red"
$ws.Range("C3").Value = "green"
$ws.Range("C4").Value = "green"
$ws.Range("C6").Value = "unknown"
$ws.Range("C7").Value = "unknown"
$ws.Range("C8").Value = "red"

$ws.Range("A5").Value = "Gulur"
$ws.Range("B5").Value = "yellow"
$ws.Range("C5").Value = "red"

$ws.Range("D1").Value = "clinical_assessment_order_simple"
$ws.Range("D2").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 1

# Column widths to fit the new headers/content (nearest values the
# engine's character-width quantization can represent).
$ws.Columns.Item(1).ColumnWidth = 27.5
$ws.Columns.Item(2).ColumnWidth = 30
$ws.Columns.Item(3).ColumnWidth = 31.5
$ws.Columns.Item(4).ColumnWidth = 32.166666666666664

# Leave the selection parked one row below the data, same as the sheet
# was left in the source workbook.
$ws.Range("D9").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. "lsh_unit_categories" gets its "order_simple" column (E) a proper
#    width, and column D widened to fit; selection parked back at D1.
# ---------------------------------------------------------------------
$unitCategories = $wb.Worksheets.Item("lsh_unit_categories")
$unitCategories.Columns.Item(4).ColumnWidth = 29
$unitCategories.Columns.Item(5).ColumnWidth = 24.666666666666668
$unitCategories.Range("D1").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. Leave the newly inserted sheet as the active/selected tab (it is
#    now the sheet that sits in front of "lsh_sheet_names").
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("D9").Select() | Out-Null
